# A new daily price quote was inserted as row 200 (pushing the existing
# rows 200..296 down to 201..297), growing the used range from A1:T296
# to A1:T297.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 200..296 down by one to make room for the new record.
$ws.Rows.Item(200).Insert()

# Populate the newly inserted row 200 with the new record's data. Columns
# A, B, C, E, F, G, H, I, J are constant across every data row in this
# sheet, so reuse those same values here.
$ws.Range("A200").Value = 5
$ws.Range("B200").Value = "Macroferia Regional de Talca"
$ws.Range("C200").Value = "Maule"
$ws.Range("D200").Value = 44846
$ws.Range("E200").Value = 7
$ws.Range("F200").Value = "Fruta"
$ws.Range("G200").Value = 100108
$ws.Range("H200").Value = "Tropicales y subtropicales"
$ws.Range("I200").Value = 100108005
$ws.Range("J200").Value = "Piña"
$ws.Range("K200").Value = "Caramelo"
$ws.Range("L200").Value = "Segunda"
$ws.Range("M200").Value = 300
$ws.Range("N200").Value = 20000
$ws.Range("O200").Value = 20000
$ws.Range("P200").Value = 20000
$ws.Range("Q200").Value = "$/caja 14 unidades"
$ws.Range("R200").Value = "Ecuador"
$ws.Range("S200").Value = 1429
$ws.Range("T200").Value = 14
